# Apply updated data dictionary edits to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (mSensPValue row): update description and missing-value designation so the
# p-value is described as coming from either an LME or ANOVA model, and note it is
# still available (not NA) for ANOVA models.
$ws.Range("B7").Value = "P-value of maternal sensitivity effect extracted from LME or ANOVA model summary"
$ws.Range("C7").Value = "NA if model did not converge/had singular fit"

# Row 13 (agePValue row): same kind of update for the age effect p-value.
$ws.Range("B13").Value = "P-value of age effect extracted from LME or ANOVA model summary"
$ws.Range("C13").Value = "NA if model did not converge/had singular fit"

# Update the active view: scroll position and selected cell
$ws.Activate()
$ws.Range("B15").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
